$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.004963351638465073, 0.004834212512253471, 0.9906736939747048, 0.9634655839616465, 0.2594227710945805)
    3  = @(0.003770732141401472, 0.00361547136886158,  0.9929146693285493, 0.9726752168337587, 0.4372503071926935)
    4  = @(0.003429803662143611, 0.003264441653402447, 0.9935552878583357, 0.9753268970896201, 0.5000027187813955)
    5  = @(0.003292163988191401, 0.003103780506519008, 0.9938139233952664, 0.9765406859767927, 0.5289152812660634)
    6  = @(0.003235555690129871, 0.003018017049313306, 0.9939203066067568, 0.977188336127628,  0.5464952269244026)
    7  = @(0.003215636986979346, 0.002978544582061501, 0.9939577344190823, 0.9774861471260623, 0.5566465993413788)
    8  = @(0.003208329051028374, 0.002958976679893763, 0.9939714694491053, 0.9776332370965403, 0.5616437937463485)
    9  = @(0.003204016641743583, 0.0029553215929497,   0.9939795698876578, 0.9776607381133388, 0.5760438850718214)
    10 = @(0.00319974991137509,  0.002959257475508847, 0.9939875851829847, 0.9776308360184502, 0.5876316320306215)
    11 = @(0.003195457044718022, 0.002964026638212515, 0.9939956362676579, 0.9775943548661169, 0.5979775984913919)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
}
